# "Generate Report for Handoff"
# The localization report now reflects that b.md has been handed off for
# localization: status moves from "Handed back: in sync with en-US" to
# "Ready for handoff", a new handoff xliff + timestamp is recorded, the
# "Content Duplicate" flag flips to False, and an out-of-date handback
# warning is recorded in the Error Detail column (with its column widened
# to fit the message) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$handoffStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8bb685aed9d0c724014778381bc99ac3cc09cde/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8f77fd8965a58397da18d65a1dddf5a291e4d432/e2e/b.md."

# --- Overview sheet: row 3 is the b.md entry ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $handoffStatus
$wsOverview.Range("F3").Value = $handoffStatus
$wsOverview.Range("G3").Value = "2016-09-06 06:44:11"

# --- zh-cn sheet: row 3 is the b.md entry ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $handoffStatus
# A leading apostrophe forces the literal text "False" to be stored as a
# string instead of being auto-coerced to a Boolean cell; ClearFormats
# removes the "quote prefix" cell style that the apostrophe entry leaves
# behind so the cell keeps using the default style, same as the rest of
# the column.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").ClearFormats()
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-06 06:43:58"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666667

# --- de-de sheet: row 3 is the b.md entry ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $handoffStatus
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").ClearFormats()
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-06 06:44:11"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666667
